$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date + Count values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-12-07T15:05:41+00:00"
$meta.Range("B21").Value = "4"

# --- Concepts sheet: insert a new "Fetal" concept as the first row,
#     shifting the existing rows down by one. ---
$concepts = $wb.Worksheets.Item("Concepts")

# Insert a blank row above the first concept row (row 2). xlShiftDown = -4121.
$concepts.Rows.Item(2).Insert(-4121, 0)

# Clone the (now shifted-down) row 3 into the newly blank row 2, so the new
# row inherits the exact same formatting as the other concept rows, then
# overwrite its Code/Display with the new "Fetal" concept.
$concepts.Range("A3:D3").Copy($concepts.Range("A2:D2"))
$concepts.Range("B2").Value = "Fetal"
$concepts.Range("C2").Value = "Fetal"
